$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 9192
$ws1.Range("F5").Value = 511
$ws1.Range("F6").Value = 719
$ws1.Range("F7").Value = 1393
$ws1.Range("F8").Value = 211
$ws1.Range("F9").Value = 72
$ws1.Range("F10").Value = 106
$ws1.Range("F11").Value = 5987
$ws1.Range("F15").Value = 4718
$ws1.Range("F21").Value = 36
$ws1.Range("F23").Value = 269
$ws1.Range("F24").Value = 19
$ws1.Range("F25").Value = 3144

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 51

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 9192
$ws4.Range("F5").Value = 51
$ws4.Range("F6").Value = 511
$ws4.Range("F7").Value = 719
$ws4.Range("F8").Value = 1393
$ws4.Range("F9").Value = 211
$ws4.Range("F10").Value = 72
$ws4.Range("F11").Value = 106
$ws4.Range("F12").Value = 5987
$ws4.Range("F16").Value = 4718
$ws4.Range("F22").Value = 36
$ws4.Range("F24").Value = 269
$ws4.Range("F25").Value = 19
$ws4.Range("F26").Value = 3146
